# Regenerate save_data to use K (strikeouts) instead of Strike# for column G,
# writing the recalculated s_vals for each start's K column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value (column G) as recomputed by the regen script.
$kValues = [ordered]@{
    2  = 1
    3  = 2
    5  = 0
    6  = 1
    7  = 0
    8  = 2
    9  = 1
    10 = 2
    11 = 2
    12 = 1
    13 = 2
    14 = 2
    15 = 0
    17 = 0
    18 = 0
    19 = 3
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 4
    26 = 1
    27 = 1
    28 = 2
    29 = 2
    30 = 0
    32 = 0
    33 = 2
    34 = 1
    35 = 0
    36 = 0
    38 = 1
    39 = 1
    40 = 4
    41 = 1
    42 = 2
    43 = 0
    44 = 2
    45 = 1
    46 = 1
    47 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
